# Apply updated conflict counts to the timetable conflict table.
# Rows 2-5: columns B,C,D,E increase 10 -> 11; columns L and P increase 5 -> 6
# Rows 12,16: columns B,C,D,E increase 5 -> 6; columns L and P increase 5 -> 6

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsTen = @(2, 3, 4, 5)
foreach ($r in $rowsTen) {
    $ws.Range("B${r}:E${r}").Value = 11
    $ws.Range("L${r}").Value = 6
    $ws.Range("P${r}").Value = 6
}

$rowsFive = @(12, 16)
foreach ($r in $rowsFive) {
    $ws.Range("B${r}:E${r}").Value = 6
    $ws.Range("L${r}").Value = 6
    $ws.Range("P${r}").Value = 6
}
